# Added Configurable zero_before_threshold parameter to enable setting dims
# before noise_threshold or First Rise Point to 0.
# This shifts the First_Noticeable_Increase_Index (col C) and its corresponding
# cumulative value (col E) for each signal segment, which in turn changes the
# derived Pulse_Width (col G = Point_Exceeds_Index - First_Noticeable_Increase_Index)
# on every Step3_DataPts_* sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("C2").Value = 87
$ws.Range("E2").Value = 0.02210712904845834
$ws.Range("G2").Value = 44
$ws.Range("C3").Value = 91
$ws.Range("E3").Value = 0.0288971309636744
$ws.Range("G3").Value = 44
$ws.Range("C4").Value = 88
$ws.Range("E4").Value = 0.0227047591095743
$ws.Range("G4").Value = 42
$ws.Range("C5").Value = 89
$ws.Range("E5").Value = 0.04152797746471708
$ws.Range("G5").Value = 17
$ws.Range("C6").Value = 87
$ws.Range("E6").Value = 0.005482234257371558
$ws.Range("G6").Value = 19

$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("C2").Value = 87
$ws.Range("E2").Value = 0.02210712904845834
$ws.Range("G2").Value = 60
$ws.Range("C3").Value = 91
$ws.Range("E3").Value = 0.0288971309636744
$ws.Range("G3").Value = 60
$ws.Range("C4").Value = 88
$ws.Range("E4").Value = 0.0227047591095743
$ws.Range("G4").Value = 62
$ws.Range("C5").Value = 89
$ws.Range("E5").Value = 0.04152797746471708
$ws.Range("G5").Value = 56
$ws.Range("C6").Value = 87
$ws.Range("E6").Value = 0.005482234257371558
$ws.Range("G6").Value = 60

$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("C2").Value = 87
$ws.Range("E2").Value = 0.02210712904845834
$ws.Range("G2").Value = 73
$ws.Range("C3").Value = 91
$ws.Range("E3").Value = 0.0288971309636744
$ws.Range("G3").Value = 63
$ws.Range("C4").Value = 88
$ws.Range("E4").Value = 0.0227047591095743
$ws.Range("G4").Value = 76
$ws.Range("C5").Value = 89
$ws.Range("E5").Value = 0.04152797746471708
$ws.Range("G5").Value = 62
$ws.Range("C6").Value = 87
$ws.Range("E6").Value = 0.005482234257371558
$ws.Range("G6").Value = 72

$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("C2").Value = 87
$ws.Range("E2").Value = 0.02210712904845834
$ws.Range("G2").Value = 87
$ws.Range("C3").Value = 91
$ws.Range("E3").Value = 0.0288971309636744
$ws.Range("G3").Value = 82
$ws.Range("C4").Value = 88
$ws.Range("E4").Value = 0.0227047591095743
$ws.Range("G4").Value = 99
$ws.Range("C5").Value = 89
$ws.Range("E5").Value = 0.04152797746471708
$ws.Range("G5").Value = 78
$ws.Range("C6").Value = 87
$ws.Range("E6").Value = 0.005482234257371558
$ws.Range("G6").Value = 83
